# Add a new "Sheet2" with a device/reader compatibility log (format parsing
# similar to micropos), placed after Sheet1, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add a brand new worksheet right after the existing Sheet1 and rename it.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws2.Name = "Sheet2"

# Header row (row 5)
$ws2.Range("B5").Value = "no"
$ws2.Range("C5").Value = "Device"
$ws2.Range("D5").Value = "baudarate"
$ws2.Range("E5").Value = "freq"

$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = "reader"
$ws2.Range("D6").Value = 9600

$ws2.Range("B7").Value = 2
$ws2.Range("C7").Value = "samsung a3"
$ws2.Range("D7").Value = 10000
$ws2.Range("E7").Value = 3850000

$ws2.Range("B8").Value = 3
$ws2.Range("C8").Value = "nokia lama"
$ws2.Range("D8").Value = 8800
$ws2.Range("E8").Value = 3200000

$ws2.Range("B9").Value = 4
$ws2.Range("C9").Value = "iphone 5se"
$ws2.Range("D9").Value = 12800
$ws2.Range("E9").Value = 4700000
$ws2.Range("F9").Value = 10000
$ws2.Range("G9").Value = 38500000

$ws2.Range("B10").Value = 5
$ws2.Range("C10").Value = "iphone 6plus"
$ws2.Range("D10").Value = 12800
$ws2.Range("E10").Value = 4700000
$ws2.Range("F10").Value = 10000
$ws2.Range("G10").Value = 38500000

$ws2.Range("B17").Value = 14
$ws2.Range("C17").Value = "iphone 4"
$ws2.Range("D17").Value = "freq tidak dapat"

$ws2.Range("B12").Value = 7
$ws2.Range("C12").Value = "oppo"
$ws2.Range("D12").Value = 10000
$ws2.Range("E12").Value = 3850000

$ws2.Range("B13").Value = 8
$ws2.Range("C13").Value = "huawei p8"
$ws2.Range("D13").Value = 10000
$ws2.Range("E13").Value = 3850000

$ws2.Range("B14").Value = 9
$ws2.Range("C14").Value = "redmi"
$ws2.Range("D14").Value = 10000
$ws2.Range("E14").Value = 3850000

$ws2.Range("F13").Value = "hang setelah command 00 88"

$ws2.Range("B15").Value = 12
$ws2.Range("C15").Value = "xpria z2"
$ws2.Range("D15").Value = 10000
$ws2.Range("E15").Value = 3850000

$ws2.Range("B16").Value = 13
$ws2.Range("C16").Value = "ace z500"

$ws2.Range("B11").Value = 6
$ws2.Range("C11").Value = "iphone 4 se"

$ws2.Range("D16").Value = "menunggu perbaikan chip emulator"

$ws2.Range("D11").Value = "freq tidak dapat"

# Column widths matching the target layout.
$ws2.Columns.Item(3).ColumnWidth = 26.140625
$ws2.Columns.Item(4).ColumnWidth = 33.140625

# Selection on the new sheet, and make it the active / selected tab.
$ws2.Range("F18").Select()
$ws2.Activate()
